$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New match data for columns HR, HS, HT (columns 226-228), rows 1-102
$newData = @(
    @(10247, 10256, 10261),
    @(2020, 2020, 2020),
    @(10, 11, 12),
    @(0, 1, 0),
    @(0, 1, 0),
    @(90, 93, 91),
    @(57, 34, 31),
    @(33, 59, 60),
    @(1, 1, 1),
    @(12, 15, 13),
    @(187, 197, 210),
    @(136, 126, 121),
    @(323, 323, 331),
    @(1.38, 1.56, 1.74),
    @(69, 100, 108),
    @(39, 60, 49),
    @(28, 30, 20),
    @(16, 14, 13),
    @(18, 13, 18),
    @(13, 14, 14),
    @(10, 9, 12),
    @(10, 8, 6),
    @(2, 1, 1),
    @(25, 23, 21),
    @(52, 60.9, 66.7),
    @(24.85, 23.07, 23.64),
    @(12.92, 14.04, 15.76),
    @(40, 31, 31),
    @(39, 32, 49),
    @(35, 37, 33),
    @(48, 40, 47),
    @(1.92, 1.74, 2.24),
    @(3.69, 2.86, 3.36),
    @(47.9, 55, 42.6),
    @(27.1, 35, 29.8),
    @(187.9, 188.6, 188.5),
    @(87.5, 89.3, 89.2),
    @(27.16, 27.66, 27.33),
    @(127, 139.1, 132.8),
    @(7, 4, 5),
    @(4, 5, 5),
    @(1, 3, 3),
    @(10, 10, 9),
    @(119, 110, 128),
    @(190, 205, 208),
    @(227, 249, 250),
    @(70.3, 77.09999999999999, 75.5),
    @(39, 32, 49),
    @(6, 14, 17),
    @(14, 14, 18),
    @(40, 31, 31),
    @(35, 37, 33),
    @(34, 38, 53),
    @(2, 1, 2),
    @(10, 9, 12),
    @(76.90000000000001, 64.3, 85.7),
    @(155, 151, 157),
    @(126, 116, 95),
    @(281, 267, 252),
    @(1.23, 1.3, 1.65),
    @(46, 51, 53),
    @(51, 37, 50),
    @(24, 28, 32),
    @(18, 13, 18),
    @(16, 14, 13),
    @(9, 4, 4),
    @(8, 3, 2),
    @(2, 7, 5),
    @(1, 3, 2),
    @(12, 14, 11),
    @(75, 28.6, 36.4),
    @(31.22, 66.75, 63),
    @(23.42, 19.07, 22.91),
    @(30, 31, 30),
    @(43, 36, 38),
    @(33, 25, 31),
    @(45, 43, 37),
    @(3.75, 3.07, 3.36),
    @(5, 10.75, 9.25),
    @(24.4, 25.6, 24.3),
    @(20, 9.300000000000001, 10.8),
    @(189, 188, 187.5),
    @(87.2, 85.09999999999999, 83.5),
    @(26.8, 24.24, 25.8),
    @(96.5, 67, 106.2),
    @(7, 11, 6),
    @(5, 6, 5),
    @(4, 3, 4),
    @(6, 2, 7),
    @(100, 119, 105),
    @(168, 152, 140),
    @(209, 173, 172),
    @(74.40000000000001, 64.8, 68.3),
    @(43, 36, 38),
    @(2, 6, 4),
    @(5, 7, 2),
    @(30, 31, 30),
    @(33, 25, 31),
    @(41, 47, 31),
    @(3, 3, 0),
    @(8, 3, 2),
    @(88.90000000000001, 75, 50)
)

for ($i = 0; $i -lt $newData.Count; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 226).Value = $newData[$i][0]
    $ws.Cells.Item($row, 227).Value = $newData[$i][1]
    $ws.Cells.Item($row, 228).Value = $newData[$i][2]
}
